$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Image" URL for James Hirchak's row (row 2) to the new SharePoint
# download link (replacing the old placehold.co link).
$ws.Range("D2").Value = "https://astutereview.sharepoint.com/sites/AstuteReviewInternal2/_layouts/15/download.aspx?UniqueId=07fdf07d%2D4215%2D4dd1%2D81cd%2Da85cc1bfbddb"

# Widen column D to comfortably fit the new, longer URL.
$ws.Columns("D").ColumnWidth = 31.92

# Leave the selection where the author left it after the edit.
$ws.Range("G23").Select()
